$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared string text updates (title volume/issue number and date range) ---
$ws.Range("A8").Value = "Volume 30   Number  14"
$ws.Range("C9").Value = "Report Covering the Week  4/3/2023  Through  4/9/2023"

# --- Crime Complaints table updates (rows 14-29) ---

# Row 14
$ws.Range("F14").Formula = "'0"
$ws.Range("C14").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("H14").Value = -100
$ws.Range("N14").Value = -84.210526315789

# Row 15
$ws.Range("C15").Value = 2
$ws.Range("G14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 100
$ws.Range("F15").Value = 5
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 11
$ws.Range("J15").Value = 17
$ws.Range("K15").Value = -35.294117647058
$ws.Range("L15").Value = 37.5
$ws.Range("M15").Value = 22.222222222222
$ws.Range("N15").Value = -50

# Row 16
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 11
$ws.Range("E16").Value = -90.909090909090
$ws.Range("F16").Value = 30
$ws.Range("G16").Value = 50
$ws.Range("H16").Value = -40
$ws.Range("I16").Value = 142
$ws.Range("J16").Value = 166
$ws.Range("K16").Value = -14.457831325301
$ws.Range("L16").Value = 97.222222222222
$ws.Range("M16").Value = 23.478260869565
$ws.Range("N16").Value = -74.087591240875

# Row 17
$ws.Range("C17").Value = 15
$ws.Range("D17").Value = 20
$ws.Range("E17").Value = -25
$ws.Range("F17").Value = 59
$ws.Range("G17").Value = 66
$ws.Range("H17").Value = -10.606060606060
$ws.Range("I17").Value = 205
$ws.Range("J17").Value = 177
$ws.Range("K17").Value = 15.819209039548
$ws.Range("L17").Value = 29.746835443038
$ws.Range("M17").Value = 49.635036496350
$ws.Range("N17").Value = -4.651162790697

# Row 18
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -28.571428571428
$ws.Range("F18").Value = 24
$ws.Range("G18").Value = 27
$ws.Range("H18").Value = -11.111111111111
$ws.Range("I18").Value = 104
$ws.Range("J18").Value = 82
$ws.Range("K18").Value = 26.829268292682
$ws.Range("L18").Value = 82.456140350877
$ws.Range("M18").Value = 5.050505050505
$ws.Range("N18").Value = -77.391304347826

# Row 19
$ws.Range("C19").Value = 19
$ws.Range("D19").Value = 18
$ws.Range("E19").Value = 5.555555555555
$ws.Range("F19").Value = 75
$ws.Range("G19").Value = 81
$ws.Range("H19").Value = -7.407407407407
$ws.Range("I19").Value = 253
$ws.Range("J19").Value = 277
$ws.Range("K19").Value = -8.664259927797
$ws.Range("L19").Value = 66.447368421052
$ws.Range("M19").Value = 91.666666666666
$ws.Range("N19").Value = 25.247524752475

# Row 20
$ws.Range("C20").Value = 13
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 225
$ws.Range("F20").Value = 43
$ws.Range("G20").Value = 26
$ws.Range("H20").Value = 65.384615384615
$ws.Range("I20").Value = 171
$ws.Range("J20").Value = 157
$ws.Range("K20").Value = 8.917197452229
$ws.Range("L20").Value = 128
$ws.Range("M20").Value = 216.666666666667
$ws.Range("N20").Value = -67.674858223062

# Row 21
$ws.Range("C21").Value = 55
$ws.Range("D21").Value = 61
$ws.Range("E21").Value = -9.836065573770
$ws.Range("F21").Value = 236
$ws.Range("G21").Value = 257
$ws.Range("H21").Value = -8.171206225680
$ws.Range("I21").Value = 889
$ws.Range("J21").Value = 880
$ws.Range("K21").Value = 1.022727272727
$ws.Range("L21").Value = 68.690702087286
$ws.Range("M21").Value = 62.522851919561
$ws.Range("N21").Value = -55.438596491228

# Row 22
$ws.Range("C22").Formula = "'0"
$ws.Range("C14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("F22").Value = 1
$ws.Range("L22").Value = -20
$ws.Range("M22").Value = -33.333333333333

# Row 23
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 9
$ws.Range("E23").Value = -66.666666666666
$ws.Range("F23").Value = 18
$ws.Range("G23").Value = 28
$ws.Range("H23").Value = -35.714285714285
$ws.Range("I23").Value = 79
$ws.Range("J23").Value = 91
$ws.Range("K23").Value = -13.186813186813
$ws.Range("L23").Value = 29.508196721311
$ws.Range("M23").Value = 46.296296296296

# Row 24
$ws.Range("C24").Value = 31
$ws.Range("D24").Value = 39
$ws.Range("E24").Value = -20.512820512820
$ws.Range("F24").Value = 158
$ws.Range("G24").Value = 179
$ws.Range("H24").Value = -11.731843575419
$ws.Range("I24").Value = 550
$ws.Range("J24").Value = 475
$ws.Range("K24").Value = 15.789473684210
$ws.Range("L24").Value = 79.153094462540
$ws.Range("M24").Value = 51.098901098901

# Row 25
$ws.Range("C25").Value = 24
$ws.Range("D25").Value = 28
$ws.Range("E25").Value = -14.285714285714
$ws.Range("F25").Value = 83
$ws.Range("G25").Value = 109
$ws.Range("H25").Value = -23.853211009174
$ws.Range("I25").Value = 286
$ws.Range("J25").Value = 295
$ws.Range("K25").Value = -3.050847457627
$ws.Range("L25").Value = 25.991189427312
$ws.Range("M25").Value = -27.040816326530

# Row 26
$ws.Range("C26").Value = 3
$ws.Range("G14").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 200
$ws.Range("F26").Value = 6
$ws.Range("G26").Value = 7
$ws.Range("H26").Value = -14.285714285714
$ws.Range("I26").Value = 16
$ws.Range("J26").Value = 26
$ws.Range("K26").Value = -38.461538461538
$ws.Range("L26").Value = 33.333333333333

# Row 27
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -33.333333333333
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -28.571428571428
$ws.Range("I27").Value = 31
$ws.Range("J27").Value = 19
$ws.Range("K27").Value = 63.157894736842
$ws.Range("L27").Value = 29.166666666666

# Row 28
$ws.Range("C28").Value = 1
$ws.Range("G14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = -50
$ws.Range("I28").Value = 6
$ws.Range("K28").Value = -50
$ws.Range("L28").Value = -53.846153846153
$ws.Range("M28").Value = -40
$ws.Range("N28").Value = -83.333333333333

# Row 29
$ws.Range("C29").Value = 1
$ws.Range("G14").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("G29").Value = 4
$ws.Range("H29").Value = -50
$ws.Range("I29").Value = 6
$ws.Range("K29").Value = -50
$ws.Range("L29").Value = -53.846153846153
$ws.Range("M29").Value = -25
$ws.Range("N29").Value = -81.818181818181
